$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.685.10"
$ws.Range("E2").Value = "  +1.06%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.644.10"
$ws.Range("E3").Value = "  +0.26%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.14"
$ws.Range("E5").Value = "  +0.79%  "
$ws.Range("E6").Value = "  -0.43%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.30"
$ws.Range("E8").Value = "  +1.02%  "
$ws.Range("E9").Value = "  +1.01%  "
$ws.Range("E10").Value = "  +0.39%  "
$ws.Range("E11").Value = "  +0.56%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.877.78"
$ws.Range("E12").Value = "  +0.37%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.606.91"
$ws.Range("E13").Value = "  -2.74%  "
$ws.Range("E14").Value = "  +0.64%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.562"
$ws.Range("E15").Value = "  +0.67%  "
$ws.Range("E16").Value = "  +0.60%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.687.04"
$ws.Range("E17").Value = "  +1.17%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "231.35"
$ws.Range("E18").Value = "  +0.59%  "
$ws.Range("E19").Value = "  +1.14%  "
$ws.Range("E20").Value = "  +0.64%  "
$ws.Range("E22").Value = "  -0.62%  "
$ws.Range("E23").Value = "  +5.29%  "
$ws.Range("E24").Value = "  -3.24%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "150.08"
$ws.Range("E25").Value = "  +1.74%  "
$ws.Range("E26").Value = "  -0.25%  "
$ws.Range("E27").Value = "  -1.15%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.67"
$ws.Range("E28").Value = "  +0.81%  "
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("E30").Value = "  +0.68%  "
$ws.Range("E31").Value = "  +0.70%  "
$ws.Range("E32").Value = "  +1.03%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.455.59"
$ws.Range("E33").Value = "  +3.28%  "
$ws.Range("E34").Value = "  +0.67%  "
$ws.Range("E35").Value = "  +0.31%  "
$ws.Range("E36").Value = "  -0.93%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.569"
$ws.Range("E37").Value = "  +1.01%  "
$ws.Range("E38").Value = "  +0.39%  "
$ws.Range("E39").Value = "  +0.56%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.894"
$ws.Range("E40").Value = "  +12.71%  "
$ws.Range("E41").Value = "  +10.24%  "
$ws.Range("E42").Value = "  +0.34%  "
$ws.Range("E43").Value = "  +0.08%  "
$ws.Range("E44").Value = "  +2.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.48"
$ws.Range("E45").Value = "  +0.96%  "
$ws.Range("E46").Value = "  +0.63%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.787.28"
$ws.Range("E47").Value = "  +0.39%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.72"
$ws.Range("E48").Value = "  +4.94%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "86.06"
$ws.Range("E49").Value = "  -1.51%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0991"
$ws.Range("E50").Value = "  +0.11%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.78"
$ws.Range("E51").Value = "  +1.35%  "
